$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("G3").Value = 'Ballari (Bellary)'
$ws.Range("G4").Value = 'Tumakuru (Tumkur)'
$ws.Range("G6").Value = 'Tumakuru (Tumkur)'
$ws.Range("G7").Value = 'Tumakuru (Tumkur)'
$ws.Range("G8").Value = 'Tumakuru (Tumkur)'
$ws.Range("G9").Value = 'Chitradurga'
$ws.Range("G10").Value = 'Chitradurga'
$ws.Range("G13").Value = 'Ballari (Bellary)'
$ws.Range("G14").Value = 'Tumakuru (Tumkur)'
$ws.Range("G16").Value = 'Ballari (Bellary)'
$ws.Range("G18").Value = 'Tumakuru (Tumkur)'
$ws.Range("G19").Value = 'Uttara Kannada (Karwar)'
$ws.Range("G21").Value = 'Chitradurga'
$ws.Range("G22").Value = 'Chitradurga'
$ws.Range("G23").Value = 'Chitradurga'
$ws.Range("G24").Value = 'Ballari (Bellary)'
$ws.Range("G25").Value = 'Tumakuru (Tumkur)'
$ws.Range("G26").Value = 'Tumakuru (Tumkur)'
$ws.Range("G29").Value = 'Ballari (Bellary)'
$ws.Range("G30").Value = 'Tumakuru (Tumkur)'
$ws.Range("G31").Value = 'Tumakuru (Tumkur)'
$ws.Range("G33").Value = 'Tumakuru (Tumkur)'
$ws.Range("G34").Value = 'Chitradurga'
$ws.Range("G35").Value = 'Ballari (Bellary)'
$ws.Range("G36").Value = 'Chitradurga'
$ws.Range("G37").Value = 'Chitradurga'
$ws.Range("G38").Value = 'Chitradurga'
$ws.Range("G41").Value = 'Tumakuru (Tumkur)'
$ws.Range("G42").Value = 'Chitradurga'
$ws.Range("G43").Value = 'Chitradurga'
$ws.Range("G44").Value = 'Tumakuru (Tumkur)'
$ws.Range("G45").Value = 'Tumakuru (Tumkur)'
$ws.Range("G46").Value = 'Chitradurga'
$ws.Range("G47").Value = 'Tumakuru (Tumkur)'
